$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "27.089.06"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "  -0.95%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.644.81"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "  -1.10%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "  +0.02%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "217.84"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "  -1.00%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "  +1.04%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "  -0.04%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.256"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "  +0.34%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.0628"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "  +0.16%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "19.98"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = "  +0.54%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0845"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "  -0.50%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "1.877.13"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "  -0.93%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "1.654.57"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "  -0.49%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "4.12"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "  -2.13%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.537"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "  +0.62%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "67.40"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "  +0.59%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "27.175.58"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "  -0.54%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "0.0₃0739"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "  +0.41%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "219.29"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "  -1.50%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "  -0.22%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "6.84"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "  +0.93%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "4.43"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "  -0.29%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = "  -1.02%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "9.19"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "  -0.69%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "148.02"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "  +0.15%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "  -0.16%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "7.41"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "  +0.04%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "  -0.67%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "15.79"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "  -1.30%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "  -1.35%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = "  -1.18%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.36"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "  -1.35%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "3.04"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "  +1.05%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "1.58"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "  +1.07%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "1.262.14"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "  -0.04%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "  -0.07%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = "  +0.55%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.543"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "  +1.03%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.842"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "  +1.49%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "  -0.15%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "0.807"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "  -0.73%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "  +4.68%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "5.37"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "  -0.33%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "1.788.73"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "  -1.01%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "62.08"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "  +0.27%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "91.93"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = "  -0.62%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "  -0.54%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "  -0.84%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "7.70"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "  +0.14%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("B50")
$r.NumberFormat = "@"
$r.Value = "BabyDogeCoin"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("C50")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.0₇0986"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "  +4.09%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("B51")
$r.NumberFormat = "@"
$r.Value = "Algorand"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("C51")
$r.NumberFormat = "@"
$r.Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0972"
$r.NumberFormat = "General"
$r.Style = "Normal"

$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "  -0.92%  "
$r.NumberFormat = "General"
$r.Style = "Normal"

